$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Fix cells whose type changes (text <-> number), using copy/paste-special to keep original style index ---
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = "0"

$ws.Range("N22").Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("N22").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = "***.*"

$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 2

$ws.Range("D15").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("D15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = "0"

$ws.Range("N22").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("N22").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"

$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = "0"

$ws.Range("N22").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("N22").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "***.*"

# --- Plain value updates (no type/style change needed) ---
$ws.Range("M14").Value = -70
$ws.Range("N14").Value = -92.105263157894
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = -33.333333333333
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 48
$ws.Range("J15").Value = 38
$ws.Range("K15").Value = 26.315789473684
$ws.Range("L15").Value = 65.51724137931
$ws.Range("M15").Value = 41.176470588235
$ws.Range("N15").Value = -27.272727272727
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 23
$ws.Range("E16").Value = -43.478260869565
$ws.Range("F16").Value = 47
$ws.Range("G16").Value = 57
$ws.Range("H16").Value = -17.543859649122
$ws.Range("I16").Value = 377
$ws.Range("J16").Value = 390
$ws.Range("K16").Value = -3.333333333333
$ws.Range("L16").Value = -17.324561403508
$ws.Range("M16").Value = -8.495145631067
$ws.Range("N16").Value = -77.275467148884
$ws.Range("C17").Value = 20
$ws.Range("D17").Value = 31
$ws.Range("E17").Value = -35.483870967741
$ws.Range("F17").Value = 81
$ws.Range("G17").Value = 112
$ws.Range("H17").Value = -27.678571428571
$ws.Range("I17").Value = 670
$ws.Range("J17").Value = 635
$ws.Range("K17").Value = 5.511811023622
$ws.Range("L17").Value = 4.037267080745
$ws.Range("M17").Value = 56.908665105386
$ws.Range("N17").Value = -21.269095182138
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -63.636363636363
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 39
$ws.Range("H18").Value = -51.282051282051
$ws.Range("I18").Value = 155
$ws.Range("J18").Value = 230
$ws.Range("K18").Value = -32.608695652173
$ws.Range("L18").Value = -32.017543859649
$ws.Range("M18").Value = -40.839694656488
$ws.Range("N18").Value = -85.110470701248
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -5.263157894736
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -11.904761904761
$ws.Range("I19").Value = 526
$ws.Range("J19").Value = 640
$ws.Range("K19").Value = -17.8125
$ws.Range("L19").Value = -31.331592689295
$ws.Range("M19").Value = 33.502538071066
$ws.Range("N19").Value = 2.935420743639
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = 58.333333333333
$ws.Range("F20").Value = 55
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 355
$ws.Range("J20").Value = 311
$ws.Range("K20").Value = 14.147909967845
$ws.Range("L20").Value = 25.441696113074
$ws.Range("M20").Value = 104.022988505747
$ws.Range("N20").Value = -76.903057905009
$ws.Range("C21").Value = 76
$ws.Range("D21").Value = 99
$ws.Range("E21").Value = -23.232323232323
$ws.Range("F21").Value = 284
$ws.Range("G21").Value = 350
$ws.Range("H21").Value = -18.857142857142
$ws.Range("I21").Value = 2137
$ws.Range("J21").Value = 2257
$ws.Range("K21").Value = -5.316792202038
$ws.Range("L21").Value = -11.584609019445
$ws.Range("M21").Value = 24.027858386535
$ws.Range("N21").Value = -62.776519770074
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 18
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -64.705882352941
$ws.Range("M22").Value = -52.631578947368
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = -4.166666666666
$ws.Range("I23").Value = 165
$ws.Range("J23").Value = 206
$ws.Range("K23").Value = -19.902912621359
$ws.Range("L23").Value = -18.71921182266
$ws.Range("M23").Value = 47.321428571428
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 41.025641025641
$ws.Range("F24").Value = 203
$ws.Range("G24").Value = 208
$ws.Range("H24").Value = -2.403846153846
$ws.Range("I24").Value = 1349
$ws.Range("J24").Value = 1377
$ws.Range("K24").Value = -2.033405954974
$ws.Range("L24").Value = -19.558735837805
$ws.Range("M24").Value = 58.333333333333
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -24
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 96
$ws.Range("H25").Value = -21.875
$ws.Range("I25").Value = 527
$ws.Range("J25").Value = 520
$ws.Range("K25").Value = 1.346153846153
$ws.Range("L25").Value = -30.474934036939
$ws.Range("C26").Value = 46
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 43.75
$ws.Range("F26").Value = 125
$ws.Range("G26").Value = 116
$ws.Range("H26").Value = 7.758620689655
$ws.Range("I26").Value = 887
$ws.Range("J26").Value = 735
$ws.Range("K26").Value = 20.680272108843
$ws.Range("L26").Value = 15.645371577575
$ws.Range("M26").Value = -17.87037037037
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = -9.090909090909
$ws.Range("I27").Value = 65
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = 14.035087719298
$ws.Range("L27").Value = 47.727272727272
$ws.Range("C28").Value = 6
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = -14.285714285714
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 36.363636363636
$ws.Range("I28").Value = 90
$ws.Range("J28").Value = 68
$ws.Range("K28").Value = 32.35294117647
$ws.Range("L28").Value = 28.571428571428
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 33.333333333333
$ws.Range("I29").Value = 34
$ws.Range("K29").Value = -19.047619047619
$ws.Range("L29").Value = -40.350877192982
$ws.Range("M29").Value = -38.181818181818
$ws.Range("N29").Value = -84.615384615384
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 33.333333333333
$ws.Range("I30").Value = 33
$ws.Range("K30").Value = -8.333333333333
$ws.Range("L30").Value = -32.653061224489
$ws.Range("M30").Value = -28.260869565217
$ws.Range("N30").Value = -83.333333333333
